# Update version string throughout the workbook from the old build
# ("mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)")
# to the new release build
# ("Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)")

$wb = $excel.ActiveWorkbook

$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

# A2: "Version: ..."
$aboutSheet.Range("A2").Value = "Version: " + $newVersion

# A6: Recommended citation text
$newCitation = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Zhaozhuang Coal Mine, China, M0430, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"
$aboutSheet.Range("A6").Value = $newCitation

# S2:S12 on the data sheet hold the build_version value for each row
for ($r = 2; $r -le 12; $r++) {
    $dataSheet.Cells.Item($r, 19).Value = $newVersion
}
